# :sparkles: Implement create broadcast users
# Rename the "whatsapp" / "email" column headers to their capitalized
# forms "Whatsapp" / "Email", and move the active selection to D1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Whatsapp"
$ws.Range("D1").Value = "Email"

[void]$ws.Range("D1").Select()
